$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the daily progress values for day "4" (column J), rows 10-13
$ws.Range("J10").Value = 2
$ws.Range("J11").Value = 3
$ws.Range("J12").Value = 1
$ws.Range("J13").Value = 4

# Update the active selection to a single cell, as recorded in the saved file
$ws.Range("D14").Select()
